$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 125: becomes a "les uitval" (class cancelled) marker row, same visual
#     pattern as row 6 (Donderdag / les uitval). Copy formatting from row 6
#     then set the text for J125.
$ws.Range("A6:J6").Copy()
$ws.Range("A125:J125").PasteSpecial(-4122)
$ws.Range("J125").Value = "les uitval"

# --- Row 126: fill in hours (2 per day) for Woensdag ---
$ws.Range("B126").Value = 2
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 2
$ws.Range("E126").Value = 2
$ws.Range("F126").Value = 2
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 2
$ws.Range("I126").Value = 2

# --- Row 127: fill in hours for Donderdag, with F127 highlighted (partial day) ---
$ws.Range("B127").Value = 4
$ws.Range("C127").Value = 4
$ws.Range("D127").Value = 4
$ws.Range("E127").Value = 4
$ws.Range("F55").Copy()
$ws.Range("F127").PasteSpecial(-4122)
$ws.Range("F127").Value = 2
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 4
$ws.Range("I127").Value = 4

# --- Row 131: fill in hours for Maandag (week "Week 41"), D131 highlighted (partial day) ---
$ws.Range("B131").Value = 4
$ws.Range("C131").Value = 4
$ws.Range("D124").Copy()
$ws.Range("D131").PasteSpecial(-4122)
$ws.Range("D131").Value = 3
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 4
$ws.Range("G131").Value = 4
$ws.Range("H131").Value = 4
$ws.Range("I131").Value = 4

# --- Update the view: active cell / selection moved to L126 ---
[void]$ws.Range("L126").Select()
